$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.2517957102541956;  C = 1.386019378421211;  D = 7.589302934471694;  E = 2.754868950507754;  F = 2.80499345858655;   G = 23 }
    3  = @{ B = 0.5243254742780024;  C = 1.822206725996549;  D = 13.98149530182482;  E = 3.739183774813004;  F = 3.789363045952605;   G = 22 }
    4  = @{ B = -0.8019598128076217; C = 1.305564707247611;  D = 6.943070400165851;  E = 2.634970664004791;  F = 2.571950116496777;   G = 21 }
    5  = @{ B = -0.02187051981424342;C = 0.4503981798184682; D = 0.5039077891406407; E = 0.709864627334424;  F = 0.7279599970962077;  G = 20 }
    6  = @{ B = 0.0766489273505621;  C = 0.6798153727021657; D = 1.082404327345544;  E = 1.040386623974734;  F = 1.065990830592977;   G = 19 }
    7  = @{ B = -0.1907918648194715; C = 0.5798188246083275; D = 0.8145864660070502; E = 0.9025444399070055; F = 0.9077227079139136;  G = 18 }
    8  = @{ B = -0.1371052455939829; C = 0.4571923072893577; D = 0.4943098393230912; E = 0.7030717170553024; F = 0.7107963775852487;  G = 17 }
    9  = @{ B = 0.1145279375814358;  C = 0.4251814381380697; D = 0.3204394163604665; E = 0.5660736845680662; F = 0.5725477730274827;  G = 16 }
    10 = @{ B = -0.01358233361013406;C = 0.3321431173215987; D = 0.2545945103064831; E = 0.5045735925576001; F = 0.5220940294212741;  G = 15 }
    11 = @{ B = -0.02420063924481428;C = 0.3457558232543146; D = 0.2682231761632309; E = 0.517902670550395;  F = 0.5384613548423456 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
